# "Generate Report for Handoff"
#
# The localization-status report has moved from "In Translation" to
# "Ready for handoff": refresh the status text and the two timestamps
# (HO Xliff generate date / zh-cn handoff date) that were stamped when the
# handoff package was produced, then re-fit the "status" columns that now
# have to hold the longer text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" --------------------
$wsOverview.Range("E2").Value = "Ready for handoff"   # zh-cn status
$wsOverview.Range("F2").Value = "Ready for handoff"   # de-de status
$wsZhCn.Range("C2").Value     = "Ready for handoff"   # Status column
$wsDeDe.Range("C2").Value     = "Ready for handoff"   # Status column

# --- Timestamps refreshed by the new handoff generation ------------------
# Latest HO Xliff Generate Date (Overview!G2) and the de-de sheet's Latest
# Handoff Datetime (H2) share the same stamp.
$wsOverview.Range("G2").Value = "2016-08-29 10:38:59"
$wsDeDe.Range("H2").Value     = "2016-08-29 10:38:59"

# zh-cn sheet's Latest Handoff Datetime got its own, slightly earlier stamp.
$wsZhCn.Range("H2").Value = "2016-08-29 10:38:55"

# --- Re-fit the columns that now hold "Ready for handoff" ----------------
$wsOverview.Columns.Item(5).ColumnWidth = 17.2159881591797  # E: zh-cn status
$wsOverview.Columns.Item(6).ColumnWidth = 17.2159881591797  # F: de-de status
$wsZhCn.Columns.Item(3).ColumnWidth     = 17.2159881591797  # C: Status
$wsDeDe.Columns.Item(3).ColumnWidth     = 17.2159881591797  # C: Status
